# Update "想去人数" (want-to-go count) values in column F on the
# "展览" and "全部类型" worksheets to reflect the latest scrape.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 461
$ws1.Range("F13").Value = 1090
$ws1.Range("F20").Value = 610
$ws1.Range("F21").Value = 268
$ws1.Range("F22").Value = 207
$ws1.Range("F23").Value = 2048
$ws1.Range("F30").Value = 2812
$ws1.Range("F33").Value = 114
$ws1.Range("F34").Value = 652
$ws1.Range("F36").Value = 1798
$ws1.Range("F38").Value = 1809
$ws1.Range("F41").Value = 837
$ws1.Range("F45").Value = 1005
$ws1.Range("F46").Value = 63
$ws1.Range("F48").Value = 3331

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 461
$ws4.Range("F14").Value = 1090
$ws4.Range("F22").Value = 610
$ws4.Range("F23").Value = 268
$ws4.Range("F24").Value = 207
$ws4.Range("F25").Value = 2048
$ws4.Range("F29").Value = 2812
$ws4.Range("F32").Value = 114
$ws4.Range("F35").Value = 652
$ws4.Range("F37").Value = 1798
$ws4.Range("F40").Value = 1809
$ws4.Range("F41").Value = 837
$ws4.Range("F44").Value = 1005
$ws4.Range("F48").Value = 3331
